$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5613.1577
$ws.Range("I74").Value = 3715.1428
$ws.Range("K74").Value = 3715.1428
$ws.Range("M74").Value = -2779.1428
$ws.Range("H77").Value = 5613.1577
$ws.Range("I77").Value = 3715.1428
$ws.Range("K77").Value = 18575.714
$ws.Range("M77").Value = -13895.714
$ws.Range("H112").Value = 8764.406000000001
$ws.Range("I112").Value = 1423.0769
$ws.Range("J112").Value = 13787.421
$ws.Range("K112").Value = 4269.2307
$ws.Range("L112").Value = 41362.263
$ws.Range("M112").Value = -3161.2307
$ws.Range("N112").Value = -43578.263

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4465.3335
$ws.Range("I2").Value = 5550.1665
$ws.Range("J2").Value = 2295.6667
$ws.Range("K2").Value = 5550.1665
$ws.Range("L2").Value = 2295.6667
$ws.Range("M2").Value = -5437.1665
$ws.Range("N2").Value = -2521.6667
$ws.Range("H32").Value = 2826670.2
$ws.Range("I32").Value = 2977295.5
$ws.Range("K32").Value = 2977295.5
$ws.Range("M32").Value = -2977008.5
$ws.Range("H45").Value = 2969.889
$ws.Range("I45").Value = 3275.7144
$ws.Range("J45").Value = 1899.5
$ws.Range("K45").Value = 3275.7144
$ws.Range("L45").Value = 1899.5
$ws.Range("M45").Value = -2898.7144
$ws.Range("N45").Value = -2653.5
$ws.Range("H61").Value = 6369.8965
$ws.Range("I61").Value = 7190.8887
$ws.Range("K61").Value = 7190.8887
$ws.Range("M61").Value = -6978.8887
$ws.Range("H74").Value = 7241.6665
$ws.Range("I74").Value = 8510.4
$ws.Range("K74").Value = 8510.4
$ws.Range("M74").Value = -7636.4
$ws.Range("H77").Value = 7241.6665
$ws.Range("I77").Value = 8510.4
$ws.Range("K77").Value = 42552
$ws.Range("M77").Value = -38184
$ws.Range("H97").Value = 1057.9412
$ws.Range("I97").Value = 1188.75
$ws.Range("J97").Value = 744
$ws.Range("K97").Value = 1188.75
$ws.Range("L97").Value = 744
$ws.Range("M97").Value = -692.75
$ws.Range("N97").Value = -1736
$ws.Range("H102").Value = 13852
$ws.Range("I102").Value = 14996.893
$ws.Range("J102").Value = 3166.3333
$ws.Range("K102").Value = 14996.893
$ws.Range("L102").Value = 3166.3333
$ws.Range("M102").Value = -13374.893
$ws.Range("N102").Value = -6410.3333
$ws.Range("H116").Value = 4465.3335
$ws.Range("I116").Value = 5550.1665
$ws.Range("J116").Value = 2295.6667
$ws.Range("K116").Value = 5550.1665
$ws.Range("L116").Value = 2295.6667
$ws.Range("M116").Value = -3256.1665
$ws.Range("N116").Value = -6883.6667
$ws.Range("H122").Value = 1874.62
$ws.Range("I122").Value = 1791.0435
$ws.Range("K122").Value = 5373.1305
$ws.Range("M122").Value = -2923.1305
$ws.Range("H132").Value = 2849.2407
$ws.Range("I132").Value = 2764.9583
$ws.Range("K132").Value = 8294.874899999999
$ws.Range("M132").Value = -5764.874899999999
$ws.Range("H136").Value = 6369.8965
$ws.Range("I136").Value = 7190.8887
$ws.Range("K136").Value = 21572.6661
$ws.Range("M136").Value = -19022.6661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4465.3335
$ws.Range("I3").Value = 5550.1665
$ws.Range("J3").Value = 2295.6667
$ws.Range("K3").Value = 5550.1665
$ws.Range("L3").Value = 2295.6667
$ws.Range("M3").Value = -5436.1665
$ws.Range("N3").Value = -2523.6667
$ws.Range("H99").Value = 34920.168
$ws.Range("I99").Value = 43296.332
$ws.Range("J99").Value = 18167.834
$ws.Range("K99").Value = 43296.332
$ws.Range("L99").Value = 18167.834
$ws.Range("M99").Value = -41798.332
$ws.Range("N99").Value = -21163.834
$ws.Range("H105").Value = 2139.7856
$ws.Range("I105").Value = 1726.8
$ws.Range("J105").Value = 3172.25
$ws.Range("K105").Value = 1726.8
$ws.Range("L105").Value = 3172.25
$ws.Range("M105").Value = 20.20000000000005
$ws.Range("N105").Value = -6666.25
$ws.Range("H107").Value = 23812264
$ws.Range("I107").Value = 3186.125
$ws.Range("J107").Value = 55557704
$ws.Range("K107").Value = 3186.125
$ws.Range("L107").Value = 55557704
$ws.Range("M107").Value = -1266.125
$ws.Range("N107").Value = -55561544

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1343
$ws.Range("I16").Value = 1181.3334
$ws.Range("K16").Value = 1181.3334
$ws.Range("M16").Value = -894.3334
$ws.Range("H22").Value = 3308
$ws.Range("I22").Value = 424.75
$ws.Range("J22").Value = 9074.5
$ws.Range("K22").Value = 424.75
$ws.Range("L22").Value = 9074.5
$ws.Range("M22").Value = -74.75
$ws.Range("N22").Value = -9774.5
$ws.Range("H31").Value = 16741.5
$ws.Range("I31").Value = 6463.9473
$ws.Range("K31").Value = 6463.9473
$ws.Range("M31").Value = -6168.9473
$ws.Range("H34").Value = 16741.5
$ws.Range("I34").Value = 6463.9473
$ws.Range("K34").Value = 6463.9473
$ws.Range("M34").Value = -6261.9473
$ws.Range("H99").Value = 11812.077
$ws.Range("I99").Value = 4409.7144
$ws.Range("K99").Value = 4409.7144
$ws.Range("M99").Value = -2911.7144
$ws.Range("H107").Value = 558.6
$ws.Range("I107").Value = 444.6154
$ws.Range("J107").Value = 1299.5
$ws.Range("K107").Value = 444.6154
$ws.Range("L107").Value = 1299.5
$ws.Range("M107").Value = 1475.3846
$ws.Range("N107").Value = -5139.5
$ws.Range("H113").Value = 1343
$ws.Range("I113").Value = 1181.3334
$ws.Range("K113").Value = 1181.3334
$ws.Range("M113").Value = 988.6666
$ws.Range("H126").Value = 11812.077
$ws.Range("I126").Value = 4409.7144
$ws.Range("K126").Value = 13229.1432
$ws.Range("M126").Value = -10759.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1240.25
$ws.Range("J55").Value = 1166.6666
$ws.Range("L55").Value = 3499.9998
$ws.Range("N55").Value = -3853.9998
$ws.Range("H86").Value = 689
$ws.Range("J86").Value = 975
$ws.Range("L86").Value = 2925
$ws.Range("N86").Value = -5297
$ws.Range("H89").Value = 689
$ws.Range("J89").Value = 975
$ws.Range("L89").Value = 8775
$ws.Range("N89").Value = -20631
$ws.Range("H97").Value = 474.25
$ws.Range("J97").Value = 474.5
$ws.Range("L97").Value = 1423.5
$ws.Range("N97").Value = -2415.5
$ws.Range("H123").Value = 3639.9092
$ws.Range("I123").Value = 2671
$ws.Range("K123").Value = 8013
$ws.Range("M123").Value = -5563

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11730.714
$ws.Range("J70").Value = 5829.25
$ws.Range("L70").Value = 5829.25
$ws.Range("N70").Value = -6369.25
$ws.Range("H73").Value = 11730.714
$ws.Range("J73").Value = 5829.25
$ws.Range("L73").Value = 5829.25
$ws.Range("N73").Value = -7701.25
$ws.Range("H102").Value = 5454.2964
$ws.Range("I102").Value = 5999.6665
$ws.Range("K102").Value = 5999.6665
$ws.Range("M102").Value = -4377.6665
$ws.Range("H107").Value = 40909.332
$ws.Range("I107").Value = 89170.5
$ws.Range("J107").Value = 2300.4
$ws.Range("K107").Value = 89170.5
$ws.Range("L107").Value = 2300.4
$ws.Range("M107").Value = -87250.5
$ws.Range("N107").Value = -6140.4
$ws.Range("H113").Value = 2578.238
$ws.Range("J113").Value = 4505.4
$ws.Range("L113").Value = 4505.4
$ws.Range("N113").Value = -8845.4
$ws.Range("H122").Value = 9969.1
$ws.Range("I122").Value = 3527.4285
$ws.Range("J122").Value = 24999.666
$ws.Range("K122").Value = 10582.2855
$ws.Range("L122").Value = 74998.99800000001
$ws.Range("M122").Value = -8132.2855
$ws.Range("N122").Value = -79898.99800000001
$ws.Range("H126").Value = 9666.333000000001
$ws.Range("I126").Value = 11999.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 35998.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -33528.5
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 298647.5
$ws.Range("J132").Value = 22496.076
$ws.Range("L132").Value = 67488.228
$ws.Range("N132").Value = -72548.228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4454.385
$ws.Range("I40").Value = 4110.778
$ws.Range("K40").Value = 4110.778
$ws.Range("M40").Value = -3974.778
$ws.Range("H46").Value = 1057.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2228.15
$ws.Range("I122").Value = 1944.9412
$ws.Range("K122").Value = 5834.8236
$ws.Range("M122").Value = -3384.8236
$ws.Range("H126").Value = 2939.8125
$ws.Range("I126").Value = 4012.375
$ws.Range("K126").Value = 12037.125
$ws.Range("M126").Value = -9567.125
$ws.Range("H132").Value = 610066.75
$ws.Range("I132").Value = 754234.9
$ws.Range("J132").Value = 7181.8184
$ws.Range("K132").Value = 2262704.7
$ws.Range("L132").Value = 21545.4552
$ws.Range("M132").Value = -2260174.7
$ws.Range("N132").Value = -26605.4552
$ws.Range("H136").Value = 7459942.5
$ws.Range("I136").Value = 9049374
$ws.Range("J136").Value = 42594
$ws.Range("K136").Value = 27148122
$ws.Range("L136").Value = 127782
$ws.Range("M136").Value = -27145572
$ws.Range("N136").Value = -132882
